$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ASTHMA IR values (re-ran the analysis with new asthma IR, excluding other states)
$ws.Range("C2").Value = 751258.50673776877
$ws.Range("C3").Value = 141224.85510650915
$ws.Range("C4").Value = 71016.182353851284
$ws.Range("C5").Value = 539017.46927740984
$ws.Range("C6").Value = 27578.037147134797
$ws.Range("C7").Value = 133224.93828604417
$ws.Range("C8").Value = 189511.85823830997
$ws.Range("C9").Value = 222768.88532188482
$ws.Range("C10").Value = 178028.82456520526
$ws.Range("C11").Value = 145.96317919143652

# Update the selected range shown when the workbook was last saved
$ws.Range("C2:C10").Select()

# Touch page setup so the sheet carries an explicit (portrait) page setup,
# matching the resave from a newer Excel build
$ws.PageSetup.Orientation = 1
